# conversion_parameters.xlsx -- "Finishing touches and running scenarios"
# Rename the "euros" cost-unit labels to "USD" and refresh several capex /
# cost figures across the parameter sheets, following the latest run.

$wb = $excel.ActiveWorkbook

# ---- "500 bar" ---------------------------------------------------------
$ws = $wb.Worksheets.Item("500 bar")
$ws.Range("A8").Value = "Compressor capex coefficient (USD per kilograms H2 per day)"
$ws.Range("B8").Value = 42189

# ---- "LH2" --------------------------------------------------------------
$ws = $wb.Worksheets.Item("LH2")
$ws.Range("A3").Value = "Capex quadratic coefficient (USD (kg H2)-2)"
$ws.Range("A4").Value = "Capex linear coefficient (USD per kg H2)"
$ws.Range("B4").Value = 1877.8
$ws.Range("A5").Value = "Capex constant (USD)"
$ws.Range("B5").Formula = "317400000"
$ws.Range("B5").Value = 317400000

# ---- "LOHC_load" ----------------------------------------------------------
$ws = $wb.Worksheets.Item("LOHC_load")
$ws.Range("A4").Value = "Capex coefficient (USD per kilograms H2 per year)"
$ws.Range("B4").Value = 0.889
$ws.Range("A7").Value = "Carrier costs (USD per kg carrier)"
$ws.Range("B7").Value = 2.116

# ---- "LOHC_unload" --------------------------------------------------------
$ws = $wb.Worksheets.Item("LOHC_unload")
$ws.Range("A4").Value = "Capex coefficient (USD per kilograms H2 per year)"
$ws.Range("B4").Value = 2.59

# ---- "NH3_load" -----------------------------------------------------------
$ws = $wb.Worksheets.Item("NH3_load")
$ws.Range("A4").Value = "Capex coefficient (USD per annual g H2)"
$ws.Range("B4").Value = 0.797906

# ---- "NH3_unload" ---------------------------------------------------------
$ws = $wb.Worksheets.Item("NH3_unload")
$ws.Range("A4").Value = "Capex coefficient (USD per hourly g H2)"
$ws.Range("B4").Value = 18191170

# The workbook was left open on the "NH3_unload" tab when the scenarios
# finished running.
$ws.Activate()
